$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: insert two new year-pairs (2022, 2021) of columns before the
#    existing data, which starts at column C (2020 at the time of editing).
# ---------------------------------------------------------------------------
$ws.Range("C1:F1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2) Bring over number formatting / fonts / borders for the two new blocks
#    from the (now shifted) neighbouring blocks so the new cells look like
#    the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("G1:H8").Copy()
$ws.Range("C1:D8").PasteSpecial(-4122)
$ws.Range("I1:J8").Copy()
$ws.Range("E1:F8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-create the merged year headers for the two new columns pairs.
# E1:F1 (2021) should look exactly like the other shifted headers (full box
# border on each half), so merge it first and then refresh its formatting
# from a genuine neighbour.
$ws.Range("E1:F1").Merge()
$ws.Range("G1:H1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# C1:D1 (2022) is merged normally.
$ws.Range("C1:D1").Merge()

# ---------------------------------------------------------------------------
# 3) Header values
# ---------------------------------------------------------------------------
$ws.Range("C1").Value2 = 2022
$ws.Range("E1").Value2 = 2021

$ws.Range("C2").Value2 = "数量"
$ws.Range("D2").Value2 = "同比去年"
$ws.Range("E2").Value2 = "数量"
$ws.Range("F2").Value2 = "同比去年"

# ---------------------------------------------------------------------------
# 4) Data rows (2022 in C/D, 2021 in E/F) — "同比去年" (= YoY change) columns
#    carry the same kind of formula already used elsewhere in the sheet.
# ---------------------------------------------------------------------------
# Row 3 - 省 (province): no change column used in the original data either
$ws.Range("C3").Value2 = 31
$ws.Range("E3").Value2 = 31

# Row 4 - 市 (city)
$ws.Range("C4").Value2 = 316
$ws.Range("E4").Value2 = 342
$ws.Range("D4").Formula = "=C4-E4"

# Row 5 - 县 (county)
$ws.Range("C5").Value2 = 3267
$ws.Range("E5").Value2 = 3271
$ws.Range("D5").Formula = "=C5-E5"
$ws.Range("F5").Formula = "=E5-G5"

# Row 6 - 镇 (town)
$ws.Range("C6").Value2 = 41313
$ws.Range("E6").Value2 = 41613
$ws.Range("D6").Formula = "=C6-E6"
$ws.Range("F6").Formula = "=E6-G6"

# Row 7 - 村 (village)
$ws.Range("C7").Value2 = 609996
$ws.Range("E7").Value2 = 633980
$ws.Range("D7").Formula = "=C7-E7"
$ws.Range("F7").Formula = "=E7-G7"

# Row 8 - 合计 (total)
$ws.Range("C8").Formula = "=SUM(C3:C7)"
$ws.Range("E8").Formula = "=SUM(E3:E7)"
$ws.Range("D8").Formula = "=C8-E8"
$ws.Range("F8").Formula = "=SUM(F3:F7)"

# ---------------------------------------------------------------------------
# 5) Selection cosmetics (matches the author ending on cell D10)
# ---------------------------------------------------------------------------
$ws.Range("D10").Select()
